$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep the literal text (Excel otherwise parses
    # numeric-looking strings like "601.32" into a float and mangles
    # them with floating point noise). Restoring the format afterwards
    # keeps the cell on the workbook's default (unstyled) style.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "68.140.46"
$ws.Range("E2").Value = "  +0.69%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "3.803.51"
$ws.Range("E3").Value = "  +0.33%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("D4") "0.997"
$ws.Range("E4").Value = "  -0.33%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "601.32"
$ws.Range("E5").Value = "  +0.85%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "166.05"
$ws.Range("E6").Value = "  -0.49%  "

# Row 7
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.517"
$ws.Range("E8").Value = "  -0.65%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.159"
$ws.Range("E9").Value = "  -0.19%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.453"
$ws.Range("E10").Value = "  +0.93%  "

# Row 11
Set-TextValue $ws.Range("D11") "6.43"
$ws.Range("E11").Value = "  +1.33%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.0000250"
$ws.Range("E12").Value = "  -1.02%  "

# Row 13
Set-TextValue $ws.Range("D13") "35.94"
$ws.Range("E13").Value = "  -0.42%  "

# Row 14
Set-TextValue $ws.Range("D14") "4.430.89"
$ws.Range("E14").Value = "  +0.11%  "

# Row 15
Set-TextValue $ws.Range("D15") "3.797.59"
$ws.Range("E15").Value = "  -0.89%  "

# Row 16
Set-TextValue $ws.Range("D16") "68.007.97"
$ws.Range("E16").Value = "  +0.50%  "

# Row 17
Set-TextValue $ws.Range("D17") "18.42"
$ws.Range("E17").Value = "  -0.86%  "

# Row 18 (no D change)
$ws.Range("E18").Value = "  +1.81%  "

# Row 19
Set-TextValue $ws.Range("D19") "7.09"
$ws.Range("E19").Value = "  +0.31%  "

# Row 20
Set-TextValue $ws.Range("D20") "464.52"
$ws.Range("E20").Value = "  +0.97%  "

# Row 21
Set-TextValue $ws.Range("D21") "9.78"
$ws.Range("E21").Value = "  -2.00%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.699"
$ws.Range("E22").Value = "  -0.15%  "

# Row 23
Set-TextValue $ws.Range("D23") "0.0000148"
$ws.Range("E23").Value = "  -2.68%  "

# Row 24
Set-TextValue $ws.Range("D24") "82.86"
$ws.Range("E24").Value = "  -0.59%  "

# Row 25
Set-TextValue $ws.Range("D25") "12.05"
$ws.Range("E25").Value = "  +0.02%  "

# Row 26
Set-TextValue $ws.Range("D26") "2.12"
$ws.Range("E26").Value = "  +0.89%  "

# Row 27 (no D change)
$ws.Range("E27").Value = "  +0.00%  "

# Row 28
Set-TextValue $ws.Range("D28") "9.98"
$ws.Range("E28").Value = "  -0.33%  "

# Row 29
Set-TextValue $ws.Range("D29") "3.946.30"
$ws.Range("E29").Value = "  +0.17%  "

# Row 30
Set-TextValue $ws.Range("D30") "7.50"
$ws.Range("E30").Value = "  +3.72%  "

# Row 31
Set-TextValue $ws.Range("D31") "2.63"
$ws.Range("E31").Value = "  -5.10%  "

# Row 32
Set-TextValue $ws.Range("D32") "2.21"
$ws.Range("E32").Value = "  -2.06%  "

# Row 33
Set-TextValue $ws.Range("D33") "29.27"
$ws.Range("E33").Value = "  -1.17%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.997"
$ws.Range("E34").Value = "  -0.24%  "

# Row 35
Set-TextValue $ws.Range("D35") "9.01"
$ws.Range("E35").Value = "  -0.64%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.0996"
$ws.Range("E36").Value = "  -0.45%  "

# Row 37 (no D change)
$ws.Range("E37").Value = "  +0.92%  "

# Row 38
Set-TextValue $ws.Range("D38") "3.29"
$ws.Range("E38").Value = "  -1.39%  "

# Row 39
Set-TextValue $ws.Range("D39") "5.79"
$ws.Range("E39").Value = "  +0.43%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.986"
$ws.Range("E40").Value = "  -0.74%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.997"
$ws.Range("E41").Value = "  -0.22%  "

# Row 42 (no D change)
$ws.Range("E42").Value = "  -0.01%  "

# Row 43
Set-TextValue $ws.Range("D43") "47.57"
$ws.Range("E43").Value = "  -1.09%  "

# Row 44
Set-TextValue $ws.Range("D44") "43.60"
$ws.Range("E44").Value = "  -0.46%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.300"
$ws.Range("E45").Value = "  +0.98%  "

# Row 46
Set-TextValue $ws.Range("D46") "151.35"
$ws.Range("E46").Value = "  +0.83%  "

# Row 47
Set-TextValue $ws.Range("D47") "8.36"
$ws.Range("E47").Value = "  +0.88%  "

# Row 48 - now EnergySwap (was Stacks)
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D48") "27.20"
$ws.Range("E48").Value = "  +1.14%  "

# Row 49 - now Stacks (was Bittensor)
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D49") "1.87"
$ws.Range("E49").Value = "  +2.87%  "

# Row 50 - now Bittensor (was EnergySwap)
$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D50") "395.93"
$ws.Range("E50").Value = "  +1.56%  "

# Row 51
Set-TextValue $ws.Range("D51") "1.35"
$ws.Range("E51").Value = "  +7.00%  "
